$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.787.63"
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("D3").Value = "2.942.34"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'592.14"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "'147.28"
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "2.940.96"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("D10").Value = "'7.04"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("D11").Value = "'0.148"
$ws.Range("E11").Value = "  +4.97%  "
$ws.Range("D12").Value = "'0.437"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "'0.0000232"
$ws.Range("E13").Value = "  +3.85%  "
$ws.Range("D14").Value = "'32.45"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "3.429.13"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "62.793.42"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "2.933.33"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "'438.22"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").Value = "'13.39"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").Value = "'0.663"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").Value = "'6.99"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "'11.17"
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("D25").Value = "'80.73"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'11.80"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").Value = "'2.12"
$ws.Range("E27").Value = "  -1.92%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").Value = "0.0₃0978"
$ws.Range("E32").Value = "  +12.24%  "
$ws.Range("D33").Value = "'26.27"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'0.989"
$ws.Range("E36").Value = "  -2.21%  "
$ws.Range("D37").Value = "'5.60"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "'3.01"
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("D39").Value = "'49.63"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "'2.01"
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("E41").Value = "  -3.50%  "
$ws.Range("D42").Value = "'8.43"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").Value = "'39.06"
$ws.Range("E44").Value = "  -7.27%  "
$ws.Range("D45").Value = "2.701.74"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'135.03"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").Value = "'0.0336"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("D48").Value = "'355.93"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "'22.63"
$ws.Range("E51").Value = "  -3.77%  "
